# DCU0015-AdicionarPersonagemAoCombate.xlsx
# "Documentação e ajustes finos"
#
# 1) E1 exception description (C24) becomes a rich-text sentence with the
#    word "dialog" in italics: "A *dialog* foi fechada"
# 2) Post-condition 5.2 (A31) text is reworded.
# 3) The alignment of the E1-description style (C24, borderId=3/fillId=7)
#    gets an explicit reading order.
# 4) A31:D31 merged range definition is re-created (moves to the end of
#    the mergeCell list).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) C24: replace plain text with rich text, italicizing "dialog" ---
$c24 = $ws.Range("C24")
$c24.Value = "A dialog foi fechada"
# "A " = 2 chars, so "dialog" starts at position 3 and is 6 chars long
$c24.Characters(3, 6).Font.Italic = $true

# --- 2) A31: reword the post-condition text ---
$ws.Range("A31").Value = "5.2 O Sistema retorna para a execução do caso de uso [DCU0005]"

# --- 3) C24 style: give the alignment an explicit reading order ---
$c24.ReadingOrder = 0

# --- 4) Recreate the A31:D31 merge so its definition is appended last ---
$mergedRange = $ws.Range("A31:D31")
$mergedRange.UnMerge()
$mergedRange.Merge()

Write-Output "edit applied"
